$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Consigner") to hold "Consignee Code".
# This shifts D:R -> E:S, including the existing "Consigner Id" header/hyperlink cell.
$ws.Columns("D:D").Insert()

# New header + value for the inserted column.
$ws.Range("D1").Value = "Consignee Code"
$ws.Range("D2").Value = "code-1011"

# Match formatting: header takes the bold style used by its neighbours (same as E1),
# the data cell takes the plain style used by column C ("Legal Name").
$ws.Range("E1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# The hyperlink on the old "Email" cell (J2) needs to move to its new location (K2) -
# the engine does not automatically re-target hyperlinks on a column insert.
$ws.Range("K2").Copy()
$ws.Range("AZ1").PasteSpecial(-4122)

$ws.Range("J2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:suman@yopmail.com")

$ws.Range("AZ1").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("AZ1").Clear()

$excel.CutCopyMode = $false

$ws.Range("D4").Select()
